$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5817
$ws.Range("I28").Value = 998.5
$ws.Range("J28").Value = 18666.334
$ws.Range("K28").Value = 998.5
$ws.Range("L28").Value = 18666.334
$ws.Range("M28").Value = -513.5
$ws.Range("N28").Value = -19636.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 55571980
$ws.Range("I40").Value = 1462.8
$ws.Range("K40").Value = 1462.8
$ws.Range("M40").Value = -1287.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1526599.9
$ws.Range("I70").Value = 1744259.9
$ws.Range("J70").Value = 2980
$ws.Range("K70").Value = 5232779.699999999
$ws.Range("L70").Value = 8940
$ws.Range("M70").Value = -5232509.699999999
$ws.Range("N70").Value = -9480

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1526599.9
$ws.Range("I73").Value = 1744259.9
$ws.Range("J73").Value = 2980
$ws.Range("K73").Value = 5232779.699999999
$ws.Range("L73").Value = 8940
$ws.Range("M73").Value = -5231843.699999999
$ws.Range("N73").Value = -10812

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5276.1177
$ws.Range("J100").Value = 6883
$ws.Range("L100").Value = 6883
$ws.Range("N100").Value = -7965

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 3541.75
$ws.Range("I129").Value = 822
$ws.Range("K129").Value = 2466
$ws.Range("M129").Value = 2534

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1914.7307
$ws.Range("I132").Value = 1558.3636
$ws.Range("K132").Value = 4675.0908
$ws.Range("M132").Value = -2145.0908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3446.46
$ws.Range("I138").Value = 1369.9546
$ws.Range("J138").Value = 5078
$ws.Range("K138").Value = 4109.8638
$ws.Range("L138").Value = 15234
$ws.Range("M138").Value = 1030.1362
$ws.Range("N138").Value = -25514

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 136887.5
$ws.Range("J140").Value = 136887.5
$ws.Range("L140").Value = 136887.5
$ws.Range("N140").Value = -147247.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6949.364
$ws.Range("J141").Value = 6665.3335
$ws.Range("L141").Value = 19996.0005
$ws.Range("N141").Value = -30356.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11963.244
$ws.Range("I32").Value = 10167.5
$ws.Range("K32").Value = 10167.5
$ws.Range("M32").Value = -9880.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17321586
$ws.Range("I61").Value = 18923190
$ws.Range("K61").Value = 18923190
$ws.Range("M61").Value = -18922978

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4339
$ws.Range("I63").Value = 4632
$ws.Range("J63").Value = 3899.5
$ws.Range("K63").Value = 4632
$ws.Range("L63").Value = 3899.5
$ws.Range("M63").Value = -3946
$ws.Range("N63").Value = -5271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4339
$ws.Range("I66").Value = 4632
$ws.Range("J66").Value = 3899.5
$ws.Range("K66").Value = 23160
$ws.Range("L66").Value = 19497.5
$ws.Range("M66").Value = -19728
$ws.Range("N66").Value = -26361.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1741.6111
$ws.Range("I74").Value = 1638.1765
$ws.Range("K74").Value = 1638.1765
$ws.Range("M74").Value = -764.1765

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1741.6111
$ws.Range("I77").Value = 1638.1765
$ws.Range("K77").Value = 8190.8825
$ws.Range("M77").Value = -3822.8825

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2743.1904
$ws.Range("I122").Value = 2743.1904
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8229.5712
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5779.5712
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 17321586
$ws.Range("I136").Value = 18923190
$ws.Range("K136").Value = 56769570
$ws.Range("M136").Value = -56767020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7459.1
$ws.Range("I86").Value = 4032
$ws.Range("J86").Value = 9743.833000000001
$ws.Range("K86").Value = 4032
$ws.Range("L86").Value = 9743.833000000001
$ws.Range("M86").Value = -2909
$ws.Range("N86").Value = -11989.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7459.1
$ws.Range("I89").Value = 4032
$ws.Range("J89").Value = 9743.833000000001
$ws.Range("K89").Value = 20160
$ws.Range("L89").Value = 48719.165
$ws.Range("M89").Value = -14544
$ws.Range("N89").Value = -59951.165

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 333282.3
$ws.Range("J105").Value = 5793.154
$ws.Range("L105").Value = 5793.154
$ws.Range("N105").Value = -9287.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32261438
$ws.Range("I31").Value = 52634420
$ws.Range("K31").Value = 52634420
$ws.Range("M31").Value = -52634125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 32261438
$ws.Range("I34").Value = 52634420
$ws.Range("K34").Value = 52634420
$ws.Range("M34").Value = -52634218

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2893.4119
$ws.Range("I58").Value = 2091.6155
$ws.Range("K58").Value = 2091.6155
$ws.Range("M58").Value = -1888.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1037.2
$ws.Range("I132").Value = 1106.5
$ws.Range("K132").Value = 3319.5
$ws.Range("M132").Value = -789.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1914.1111
$ws.Range("I134").Value = 1563.2142
$ws.Range("K134").Value = 4689.642599999999
$ws.Range("M134").Value = -2154.642599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2893.4119
$ws.Range("I136").Value = 2091.6155
$ws.Range("K136").Value = 6274.8465
$ws.Range("M136").Value = -3724.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 9552.546
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9552.546
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 28657.638
$ws.Range("N34").Value = -28825.638
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 36.2
$ws.Range("J38").Value = 15
$ws.Range("L38").Value = 45
$ws.Range("N38").Value = -739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4794753
$ws.Range("I107").Value = 3799.25
$ws.Range("J107").Value = 6072341
$ws.Range("K107").Value = 11397.75
$ws.Range("L107").Value = 18217023
$ws.Range("M107").Value = -9477.75
$ws.Range("N107").Value = -18220863

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 27878.5
$ws.Range("J122").Value = 439
$ws.Range("L122").Value = 3951
$ws.Range("N122").Value = -8851

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 14426.4
$ws.Range("I125").Value = 9599.666999999999
$ws.Range("J125").Value = 21666.5
$ws.Range("K125").Value = 28799.001
$ws.Range("L125").Value = 64999.5
$ws.Range("M125").Value = -23879.001
$ws.Range("N125").Value = -74839.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 9447
$ws.Range("I129").Value = 4435.6
$ws.Range("J129").Value = 17799.334
$ws.Range("K129").Value = 13306.8
$ws.Range("L129").Value = 53398.00199999999
$ws.Range("M129").Value = -8306.800000000001
$ws.Range("N129").Value = -63398.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1937
$ws.Range("I132").Value = 1916
$ws.Range("K132").Value = 17244
$ws.Range("M132").Value = -14714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 17650974
$ws.Range("I140").Value = 21429846
$ws.Range("K140").Value = 64289538
$ws.Range("M140").Value = -64284358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2301.6296
$ws.Range("I102").Value = 2285.76
$ws.Range("K102").Value = 2285.76
$ws.Range("M102").Value = -663.7600000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 896.13336
$ws.Range("I46").Value = 498
$ws.Range("K46").Value = 498
$ws.Range("M46").Value = -310

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3155.6897
$ws.Range("I61").Value = 2662.739
$ws.Range("K61").Value = 2662.739
$ws.Range("M61").Value = -2460.739

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2781183.2
$ws.Range("J68").Value = 5698.5713
$ws.Range("L68").Value = 5698.5713
$ws.Range("N68").Value = -7196.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2781183.2
$ws.Range("J71").Value = 5698.5713
$ws.Range("L71").Value = 28492.8565
$ws.Range("N71").Value = -35980.85649999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3155.6897
$ws.Range("I113").Value = 2662.739
$ws.Range("K113").Value = 2662.739
$ws.Range("M113").Value = -492.739

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3492.1282
$ws.Range("I122").Value = 3494.5527
$ws.Range("K122").Value = 10483.6581
$ws.Range("M122").Value = -8033.658100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3278.3774
$ws.Range("I132").Value = 1855.1389
$ws.Range("J132").Value = 6292.294
$ws.Range("K132").Value = 5565.4167
$ws.Range("L132").Value = 18876.882
$ws.Range("M132").Value = -3035.4167
$ws.Range("N132").Value = -23936.882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4930.1875
$ws.Range("I136").Value = 2049.375
$ws.Range("K136").Value = 6148.125
$ws.Range("M136").Value = -3598.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 299999
$ws.Range("J140").Value = 299999
$ws.Range("L140").Value = 299999
$ws.Range("N140").Value = -310359

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3664.5454
$ws.Range("I122").Value = 2520.7144
$ws.Range("J122").Value = 5666.25
$ws.Range("K122").Value = 7562.1432
$ws.Range("L122").Value = 16998.75
$ws.Range("M122").Value = -5112.1432
$ws.Range("N122").Value = -21898.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 139387.67
$ws.Range("J132").Value = 557595.6
$ws.Range("L132").Value = 1672786.8
$ws.Range("N132").Value = -1677846.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 347668.75
$ws.Range("J136").Value = 1430720.1
$ws.Range("L136").Value = 4292160.300000001
$ws.Range("N136").Value = -4297260.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 88948.75
$ws.Range("J141").Value = 88948.75
$ws.Range("L141").Value = 88948.75
$ws.Range("N141").Value = -99308.75
